$wb = $excel.ActiveWorkbook

# --- Sheet "Create" ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("D2").Value = "Y1"
$wsCreate.Range("E2").Value = "Y2"
$wsCreate.Range("G2").Value = "Rest"

# --- Sheet "Edit" ---
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("D2").Value = "Y1"
$wsEdit.Range("E2").Value = "Y2"
$wsEdit.Range("G2").Value = "Rest"
$wsEdit.Range("H2").Value = "XNX"
$wsEdit.Range("I2").Value = "XNA"

# --- Sheet "Delete" ---
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("D2").Value = "XNA"
$wsDelete.Range("E2").Value = "XNX"
$wsDelete.Range("G2").Value = "Rest"

# --- Selections / active sheet ---
# "Create" is no longer the active/selected tab; its selection moves to G3
$wsCreate.Range("G3").Select()

# "Delete" becomes the active/selected tab with selection at G2
$wsDelete.Activate()
$wsDelete.Range("G2").Select()
